$wb = $excel.ActiveWorkbook

# --- Refresh the "Elapsed Duration(Hrs)" column (G) on each report sheet ---
# Every open-outage row's elapsed-duration string advanced by the same
# 0:06:55 tick between the previous export and this one.

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3883:31:01"
$ws1.Range("G3").Value = "23:03:39"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12064:54:42"
$ws2.Range("G3").Value = "3194:38:11"
$ws2.Range("G4").Value = "432:49:45"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2910:44:31"
$ws4.Range("G3").Value = "137:56:46"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "384:43:30"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "25:15:48"

# --- Sheet R1 gained a new (blank) trailing row, extending the used range
# from A1:L4 to A1:L5 ---
for ($c = 1; $c -le 12; $c++) {
    $ws1.Cells.Item(5, $c).Font.Bold = $false
}
